# Update "想去人数" (interest count) figures on the "展览" and "全部类型" sheets
# to match the latest generated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rId1 / sheet1.xml) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value  = 34
$wsExhibit.Range("F5").Value  = 11298
$wsExhibit.Range("F6").Value  = 201
$wsExhibit.Range("F7").Value  = 314
$wsExhibit.Range("F9").Value  = 11221
$wsExhibit.Range("F10").Value = 456
$wsExhibit.Range("F11").Value = 1147
$wsExhibit.Range("F12").Value = 61
$wsExhibit.Range("F13").Value = 1732
$wsExhibit.Range("F14").Value = 5608
$wsExhibit.Range("F16").Value = 3459

# --- Sheet "全部类型" (rId4 / sheet4.xml) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 34
$wsAll.Range("F7").Value  = 11298
$wsAll.Range("F8").Value  = 201
$wsAll.Range("F9").Value  = 314
$wsAll.Range("F11").Value = 11221
$wsAll.Range("F12").Value = 456
$wsAll.Range("F13").Value = 1147
$wsAll.Range("F14").Value = 61
$wsAll.Range("F15").Value = 1732
$wsAll.Range("F16").Value = 5608
$wsAll.Range("F18").Value = 3459
